$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.986.41'
$ws.Range('E2').Value = '  +0.65%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.741.83'
$ws.Range('E3').Value = '  +0.46%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.27'
$ws.Range('E5').Value = '  +4.56%  '

# Row 6
$ws.Range('E6').Value = '  -0.03%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5020'
$ws.Range('E7').Value = '  -1.77%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2741'
$ws.Range('E8').Value = '  +0.38%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06195'
$ws.Range('E9').Value = '  +1.55%  '

# Row 10
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07262'
$ws.Range('E10').Value = '  +1.35%  '

# Row 11
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.741.63'
$ws.Range('E11').Value = '  +0.40%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.6546'
$ws.Range('E12').Value = '  +3.31%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.19'
$ws.Range('E13').Value = '  +2.01%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.714'
$ws.Range('E14').Value = '  +2.99%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.76'
$ws.Range('E15').Value = '  +1.03%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.03%  '

# Row 17
$ws.Range('E17').Value = '  +0.04%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.001.51'
$ws.Range('E18').Value = '  +0.65%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.91'
$ws.Range('E19').Value = '  +1.71%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006856'
$ws.Range('E20').Value = '  +2.38%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.618'
$ws.Range('E21').Value = '  +9.01%  '

# Row 22
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.967.54'
$ws.Range('E22').Value = '  +0.19%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.764'
$ws.Range('E23').Value = '  +1.83%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.410'
$ws.Range('E24').Value = '  +3.96%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.90'
$ws.Range('E25').Value = '  -2.89%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.499'
$ws.Range('E26').Value = '  -0.60%  '

# Row 27
$ws.Range('E27').Value = '  +1.12%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.791'
$ws.Range('E28').Value = '  +2.80%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '105.50'
$ws.Range('E29').Value = '  +0.58%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.982'
$ws.Range('E30').Value = '  +1.91%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08166'
$ws.Range('E31').Value = '  -2.01%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.704'
$ws.Range('E32').Value = '  +3.01%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04750'
$ws.Range('E33').Value = '  +4.63%  '

# Row 34
$ws.Range('E34').Value = '  +0.33%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9996'
$ws.Range('E35').Value = '  +2.22%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6143'
$ws.Range('E36').Value = '  -0.73%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.757'
$ws.Range('E37').Value = '  +2.72%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01624'
$ws.Range('E38').Value = '  +1.96%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.939'
$ws.Range('E39').Value = '  +1.54%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.000'
$ws.Range('E40').Value = '  +0.01%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '101.23'
$ws.Range('E41').Value = '  +3.69%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8012'
$ws.Range('E42').Value = '  +9.50%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3925'
$ws.Range('E43').Value = '  +2.68%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.051'
$ws.Range('E44').Value = '  +2.43%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1179'
$ws.Range('E45').Value = '  +4.94%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.402'
$ws.Range('E46').Value = '  +4.22%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.92'
$ws.Range('E47').Value = '  +2.50%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05293'
$ws.Range('E48').Value = '  +0.53%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.97'
$ws.Range('E49').Value = '  +1.94%  '

# Row 50
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3493'
$ws.Range('E50').Value = '  +2.60%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.655'
$ws.Range('E51').Value = '  +1.30%  '
